$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H75").Value = 10285
$ws.Range("I75").Value = 10285
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 10285
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = -9349
$ws.Range("N75").Value = $null

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H78").Value = 10285
$ws.Range("I78").Value = 10285
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 30855
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = -26175
$ws.Range("N78").Value = $null

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 83341150
$ws.Range("I86").Value = 76931144
$ws.Range("J86").Value = 95245460
$ws.Range("K86").Value = 76931144
$ws.Range("L86").Value = 95245460
$ws.Range("M86").Value = -76930021
$ws.Range("N86").Value = -95247706

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 500
$ws.Range("I87").Value = 500
$ws.Range("K87").Value = 500
$ws.Range("M87").Value = 748

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 83341150
$ws.Range("I89").Value = 76931144
$ws.Range("J89").Value = 95245460
$ws.Range("K89").Value = 384655720
$ws.Range("L89").Value = 476227300
$ws.Range("M89").Value = -384650104
$ws.Range("N89").Value = -476238532

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H90").Value = 500
$ws.Range("I90").Value = 500
$ws.Range("K90").Value = 1500
$ws.Range("M90").Value = 4740

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 47620788
$ws.Range("I107").Value = 2375.8
$ws.Range("J107").Value = 166666820
$ws.Range("K107").Value = 2375.8
$ws.Range("L107").Value = 166666820
$ws.Range("M107").Value = -455.8000000000002
$ws.Range("N107").Value = -166670660

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1545543.6
$ws.Range("I137").Value = 7074.6787
$ws.Range("J137").Value = 2709790.5
$ws.Range("K137").Value = 21224.0361
$ws.Range("L137").Value = 8129371.5
$ws.Range("M137").Value = -18674.0361
$ws.Range("N137").Value = -8134471.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 12046.071
$ws.Range("J138").Value = 4774.6
$ws.Range("L138").Value = 14323.8
$ws.Range("N138").Value = -24603.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2895.8696
$ws.Range("I32").Value = 1591.9818
$ws.Range("K32").Value = 1591.9818
$ws.Range("M32").Value = -1304.9818

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1251749.6
$ws.Range("I61").Value = 43489.58
$ws.Range("J61").Value = 2679693.2
$ws.Range("K61").Value = 43489.58
$ws.Range("L61").Value = 2679693.2
$ws.Range("M61").Value = -43277.58
$ws.Range("N61").Value = -2680117.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1251749.6
$ws.Range("I136").Value = 43489.58
$ws.Range("J136").Value = 2679693.2
$ws.Range("K136").Value = 130468.74
$ws.Range("L136").Value = 8039079.600000001
$ws.Range("M136").Value = -127918.74
$ws.Range("N136").Value = -8044179.600000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5665.788
$ws.Range("I31").Value = 4131.857
$ws.Range("J31").Value = 5968.2534
$ws.Range("K31").Value = 4131.857
$ws.Range("L31").Value = 5968.2534
$ws.Range("M31").Value = -3836.857
$ws.Range("N31").Value = -6558.2534

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 5665.788
$ws.Range("I34").Value = 4131.857
$ws.Range("J34").Value = 5968.2534
$ws.Range("K34").Value = 4131.857
$ws.Range("L34").Value = 5968.2534
$ws.Range("M34").Value = -3929.857
$ws.Range("N34").Value = -6372.2534

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2392.524
$ws.Range("I58").Value = 2180.6155
$ws.Range("J58").Value = 2736.875
$ws.Range("K58").Value = 2180.6155
$ws.Range("L58").Value = 2736.875
$ws.Range("M58").Value = -1977.6155
$ws.Range("N58").Value = -3142.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1210.7142
$ws.Range("I134").Value = 1100
$ws.Range("J134").Value = 1487.5
$ws.Range("K134").Value = 3300
$ws.Range("L134").Value = 4462.5
$ws.Range("M134").Value = -765
$ws.Range("N134").Value = -9532.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2392.524
$ws.Range("I136").Value = 2180.6155
$ws.Range("J136").Value = 2736.875
$ws.Range("K136").Value = 6541.8465
$ws.Range("L136").Value = 8210.625
$ws.Range("M136").Value = -3991.8465
$ws.Range("N136").Value = -13310.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 2723.5
$ws.Range("I7").Value = 298.33334
$ws.Range("J7").Value = 9999
$ws.Range("K7").Value = 895.0000200000001
$ws.Range("L7").Value = 29997
$ws.Range("M7").Value = -783.0000200000001
$ws.Range("N7").Value = -30221

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 600
$ws.Range("I92").Value = 600
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 1800
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = -552
$ws.Range("N92").Value = $null

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 8036.143
$ws.Range("I129").Value = 1494
$ws.Range("J129").Value = 12062.077
$ws.Range("K129").Value = 4482
$ws.Range("L129").Value = 36186.231
$ws.Range("M129").Value = 518
$ws.Range("N129").Value = -46186.231

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 3638565.2
$ws.Range("I131").Value = 10102388
$ws.Range("J131").Value = 2665
$ws.Range("K131").Value = 30307164
$ws.Range("L131").Value = 7995
$ws.Range("M131").Value = -30302124
$ws.Range("N131").Value = -18075

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 19293584
$ws.Range("I80").Value = 22912.285
$ws.Range("J80").Value = 41776036
$ws.Range("K80").Value = 22912.285
$ws.Range("L80").Value = 41776036
$ws.Range("M80").Value = -21914.285
$ws.Range("N80").Value = -41778032

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 19293584
$ws.Range("I83").Value = 22912.285
$ws.Range("J83").Value = 41776036
$ws.Range("K83").Value = 114561.425
$ws.Range("L83").Value = 208880180
$ws.Range("M83").Value = -109569.425
$ws.Range("N83").Value = -208890164

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 125001976
$ws.Range("I102").Value = 166667460
$ws.Range("J102").Value = 5555
$ws.Range("K102").Value = 166667460
$ws.Range("L102").Value = 5555
$ws.Range("M102").Value = -166665838
$ws.Range("N102").Value = -8799

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 4246.304
$ws.Range("I113").Value = 4223.5293
$ws.Range("J113").Value = 4310.8335
$ws.Range("K113").Value = 4223.5293
$ws.Range("L113").Value = 4310.8335
$ws.Range("M113").Value = -2053.5293
$ws.Range("N113").Value = -8650.833500000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2852445.8
$ws.Range("I132").Value = 3873.5
$ws.Range("K132").Value = 11620.5
$ws.Range("M132").Value = -9090.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H134").Value = 92499.75
$ws.Range("J134").Value = 92499.75
$ws.Range("L134").Value = 277499.25
$ws.Range("N134").Value = -282569.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 99999.5
$ws.Range("J136").Value = 99999.5
$ws.Range("L136").Value = 299998.5
$ws.Range("N136").Value = -305098.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 9872.416999999999
$ws.Range("I61").Value = 5179.5713
$ws.Range("J61").Value = 16442.4
$ws.Range("K61").Value = 5179.5713
$ws.Range("L61").Value = 16442.4
$ws.Range("M61").Value = -4977.5713
$ws.Range("N61").Value = -16846.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 9872.416999999999
$ws.Range("I113").Value = 5179.5713
$ws.Range("J113").Value = 16442.4
$ws.Range("K113").Value = 5179.5713
$ws.Range("L113").Value = 16442.4
$ws.Range("M113").Value = -3009.5713
$ws.Range("N113").Value = -20782.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3319.5
$ws.Range("I122").Value = 2933.5
$ws.Range("J122").Value = 5120.8335
$ws.Range("K122").Value = 8800.5
$ws.Range("L122").Value = 15362.5005
$ws.Range("M122").Value = -6350.5
$ws.Range("N122").Value = -20262.5005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5474.722
$ws.Range("I132").Value = 4503.5386
$ws.Range("J132").Value = 7999.8
$ws.Range("K132").Value = 13510.6158
$ws.Range("L132").Value = 23999.4
$ws.Range("M132").Value = -10980.6158
$ws.Range("N132").Value = -29059.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 795212.1
$ws.Range("I107").Value = 1792
$ws.Range("J107").Value = 1144317
$ws.Range("K107").Value = 5376
$ws.Range("L107").Value = 3432951
$ws.Range("M107").Value = -3456
$ws.Range("N107").Value = -3436791

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 82997.8
$ws.Range("J123").Value = 82997.8
$ws.Range("L123").Value = 82997.8
$ws.Range("N123").Value = -92797.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4117.421
$ws.Range("I136").Value = 3209.889
$ws.Range("J136").Value = 4934.2
$ws.Range("K136").Value = 9629.667000000001
$ws.Range("L136").Value = 14802.6
$ws.Range("M136").Value = -7079.667000000001
$ws.Range("N136").Value = -19902.6
